$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: collector (D2), number (E2)
$ws.Range("D2").Value = "LiRE m Daly"
# E2 looks numeric ("69058") - force text format first so Excel stores it
# as a string instead of auto-converting it to a number.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "69058"

# Row 3: collector (D3), number (E3)
$ws.Range("D3").Value = "preto Oliveira; A RS"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8031"

# Row 4: number (E4), addcoll (F4)
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2"
$ws.Range("F4").Value = "et al."

# Row 5: number (E5)
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "90957"
